# Update the team-specific transition matrix with probabilities recomputed
# after simulating more games (see commit message: "added more games, sped
# up simulate game logic, and drafted optimization logic").
#
# Each row of the sheet is a probability distribution (counts-of-outcome /
# total-games-in-row) over the possible next states; re-running the
# simulation with a larger sample changed the observed frequencies in the
# cells below.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 0.1818181818181818
$ws.Range("C2").Value = 0.5454545454545454
$ws.Range("P2").Value = 0.2727272727272727

$ws.Range("P3").Value = 0.8333333333333334
$ws.Range("S3").Value = 0.1666666666666667

$ws.Range("F6").Value = 0.125
$ws.Range("J6").Value = 0.1875
$ws.Range("O6").Value = 0.0625
$ws.Range("Q6").Value = 0.25
$ws.Range("R6").Value = 0.0625
$ws.Range("S6").Value = 0.3125

$ws.Range("B7").Value = 0.06666666666666667
$ws.Range("F7").Value = 0.2
$ws.Range("J7").Value = 0.06666666666666667
$ws.Range("O7").Value = 0.06666666666666667
$ws.Range("Q7").Value = 0.1333333333333333
$ws.Range("R7").Value = 0.1333333333333333
$ws.Range("S7").Value = 0.3333333333333333

$ws.Range("B8").Value = 0.03703703703703703
$ws.Range("F8").Value = 0.03703703703703703
$ws.Range("J8").Value = 0.03703703703703703
$ws.Range("O8").Value = 0.03703703703703703
$ws.Range("Q8").Value = 0.1851851851851852
$ws.Range("R8").Value = 0.1111111111111111
$ws.Range("S8").Value = 0.5555555555555556

$ws.Range("B9").Value = 0.05263157894736842
$ws.Range("D9").Value = 0.05263157894736842
$ws.Range("J9").Value = 0.05263157894736842
$ws.Range("R9").Value = 0.1578947368421053
$ws.Range("S9").Value = 0.6842105263157895

$ws.Range("B10").Value = 0.05952380952380952
$ws.Range("D10").Value = 0.0119047619047619
$ws.Range("F10").Value = 0.08333333333333333
$ws.Range("J10").Value = 0.09523809523809523
$ws.Range("O10").Value = 0.02380952380952381
$ws.Range("Q10").Value = 0.130952380952381
$ws.Range("R10").Value = 0.1428571428571428
$ws.Range("S10").Value = 0.4523809523809524

$ws.Range("G11").Value = 0.1379310344827586
$ws.Range("J11").Value = 0.1379310344827586
$ws.Range("K11").Value = 0.1724137931034483
$ws.Range("L11").Value = 0.5517241379310345

$ws.Range("G12").Value = 0.5882352941176471
$ws.Range("J12").Value = 0.3529411764705883
$ws.Range("S12").Value = 0.05882352941176471

$ws.Range("H15").Value = 0.1333333333333333
$ws.Range("J15").Value = 0.2666666666666667
$ws.Range("O15").Value = 0.06666666666666667
$ws.Range("S15").Value = 0.3333333333333333

$ws.Range("H16").Value = 0.1428571428571428
$ws.Range("I16").Value = 0.1428571428571428
$ws.Range("J16").Value = 0.5714285714285714
$ws.Range("O16").Value = 0.1428571428571428

$ws.Range("H17").Value = 0.1904761904761905
$ws.Range("I17").Value = 0.04761904761904762
$ws.Range("J17").Value = 0.2857142857142857
$ws.Range("K17").Value = 0.2380952380952381
$ws.Range("M17").Value = 0.04761904761904762
$ws.Range("O17").Value = 0.09523809523809523
$ws.Range("S17").Value = 0.09523809523809523

$ws.Range("I18").Value = 0.2
$ws.Range("J18").Value = 0.5
$ws.Range("K18").Value = 0.1
$ws.Range("O18").Value = 0.1
$ws.Range("S18").Value = 0.1

$ws.Range("F19").Value = 0.0202020202020202
$ws.Range("H19").Value = 0.2121212121212121
$ws.Range("I19").Value = 0.1313131313131313
$ws.Range("J19").Value = 0.3737373737373738
$ws.Range("K19").Value = 0.1313131313131313
$ws.Range("M19").Value = 0.0101010101010101
$ws.Range("O19").Value = 0.0101010101010101
$ws.Range("S19").Value = 0.1111111111111111
